# "Some other typos fixed #1"
# Fix a handful of "whse shpt. line" / "Line" typos in the Scenario column
# (Table2[Scenario]) of the ATDD Scenarios sheet:
#   "whse shpt. line"  -> "whse. shpt. line"
#   "...Line with ..."  -> "...line with ..."  (lower-case "line")
# These cells feed the [ATDD Format]/[Code Format] calculated columns
# (I and J), which recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E3").Value  = "Delete by user with no allowance manually created whse. shpt. line"
$ws.Range("E10").Value = "Delete by user with no allowance automatically created whse. shpt. line"
$ws.Range("E17").Value = "Delete by user with allowance manually created whse. shpt. line"
$ws.Range("E24").Value = "Delete by user with allowance automatically created whse. shpt. line"
$ws.Range("E36").Value = "Delete  by user with no allowance manually created whse. shpt. line"
$ws.Range("E43").Value = "Delete by user with no allowance automatically created whse. shpt. line"
$ws.Range("E50").Value = "Delete by user with allowance manually created whse. shpt. line"
$ws.Range("E57").Value = "Delete  by user with allowance automatically created whse. shpt. line with confirmation"
$ws.Range("E64").Value = "Delete  by user with allowance automatically created whse. shpt. line with no confirmation"

# Update the saved selection to match where the user was working
# (scrolled down to the newly-edited rows, selection on G63).
$ws.Range("G63").Select() | Out-Null
